$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the previously-missing readings for the existing row (2025-10-31)
$ws.Range("C2").Value = 96.09999999999999
$ws.Range("E2").Value = 27.5

# New daily readings for 2025-11-01 .. 2025-11-06
$data = @(
    @(45962, 95.3, 94.8, 28.6, 27.6),
    @(45963, 94.7, 95.59999999999999, 28, 28),
    @(45964, 95, 95.2, 28, 27.2),
    @(45965, 94.90000000000001, 95.09999999999999, 28.3, 27.2),
    @(45966, 94.40000000000001, 95, 28.4, 27.7),
    @(45967, 94.3, 94.90000000000001, 28.5, 26.9)
)

$r = 3
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 1).NumberFormat = $ws.Cells.Item(2, 1).NumberFormat
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $r = $r + 1
}
